$wb = $excel.ActiveWorkbook

# --- Sheet 1: ALC ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H17").Value = 494.7551
$ws.Range("J17").Value = 494.7551
$ws.Range("L17").Value = 1484.2653
$ws.Range("N17").Value = -1820.2653
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("N51").ClearContents()
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H113").Value = 3284.1667
$ws.Range("I113").Value = 2568.3333
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 2568.3333
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = 685.6667000000002
$ws.Range("N113").Value = -10508
$ws.Range("H137").Value = 1711.0625
$ws.Range("I137").Value = 989
$ws.Range("J137").Value = 1951.75
$ws.Range("K137").Value = 2967
$ws.Range("L137").Value = 5855.25
$ws.Range("M137").Value = -417
$ws.Range("N137").Value = -10955.25
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

# --- Sheet 2: ARM ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H2").Value = 13176
$ws.Range("I2").Value = 599.1667
$ws.Range("J2").Value = 50906.5
$ws.Range("K2").Value = 599.1667
$ws.Range("L2").Value = 50906.5
$ws.Range("M2").Value = -486.1667
$ws.Range("N2").Value = -51132.5
$ws.Range("H32").Value = 2603.1975
$ws.Range("I32").Value = 2532.1177
$ws.Range("J32").Value = 2975
$ws.Range("K32").Value = 2532.1177
$ws.Range("L32").Value = 2975
$ws.Range("M32").Value = -2245.1177
$ws.Range("N32").Value = -3549
$ws.Range("H45").Value = 1080.9474
$ws.Range("I45").Value = 1035.8667
$ws.Range("K45").Value = 1035.8667
$ws.Range("M45").Value = -658.8667
$ws.Range("H74").Value = 1217.7037
$ws.Range("I74").Value = 760
$ws.Range("K74").Value = 760
$ws.Range("M74").Value = 114
$ws.Range("H77").Value = 1217.7037
$ws.Range("I77").Value = 760
$ws.Range("K77").Value = 3800
$ws.Range("M77").Value = 568
$ws.Range("H110").Value = 1708.7142
$ws.Range("I110").Value = 1149.5
$ws.Range("J110").Value = 2454.3333
$ws.Range("K110").Value = 1149.5
$ws.Range("L110").Value = 2454.3333
$ws.Range("M110").Value = 895.5
$ws.Range("N110").Value = -6544.3333
$ws.Range("H116").Value = 13176
$ws.Range("I116").Value = 599.1667
$ws.Range("J116").Value = 50906.5
$ws.Range("K116").Value = 599.1667
$ws.Range("L116").Value = 50906.5
$ws.Range("M116").Value = 1694.8333
$ws.Range("N116").Value = -55494.5
$ws.Range("H122").Value = 751.871
$ws.Range("I122").Value = 755.1111
$ws.Range("K122").Value = 2265.3333
$ws.Range("M122").Value = 184.6667000000002
$ws.Range("H138").Value = 55551.668
$ws.Range("J138").Value = 55551.668
$ws.Range("L138").Value = 55551.668
$ws.Range("N138").Value = -65831.66800000001

# --- Sheet 3: BSM ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H3").Value = 13176
$ws.Range("I3").Value = 599.1667
$ws.Range("J3").Value = 50906.5
$ws.Range("K3").Value = 599.1667
$ws.Range("L3").Value = 50906.5
$ws.Range("M3").Value = -485.1667
$ws.Range("N3").Value = -51134.5
$ws.Range("H86").Value = 4356.2607
$ws.Range("I86").Value = 4359.75
$ws.Range("K86").Value = 4359.75
$ws.Range("M86").Value = -3236.75
$ws.Range("H89").Value = 4356.2607
$ws.Range("I89").Value = 4359.75
$ws.Range("K89").Value = 21798.75
$ws.Range("M89").Value = -16182.75
$ws.Range("H99").Value = 31251114
$ws.Range("J99").Value = 1137
$ws.Range("L99").Value = 1137
$ws.Range("N99").Value = -4133
$ws.Range("H105").Value = 250003060
$ws.Range("I105").Value = 250003060
$ws.Range("K105").Value = 250003060
$ws.Range("M105").Value = -250001313
$ws.Range("H107").Value = 1638.3334
$ws.Range("I107").Value = 1440.3334
$ws.Range("K107").Value = 1440.3334
$ws.Range("M107").Value = 479.6666
$ws.Range("H134").Value = 1433.6938
$ws.Range("I134").Value = 963.0769
$ws.Range("K134").Value = 2889.2307
$ws.Range("M134").Value = -354.2307000000001

# --- Sheet 4: CRP ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H99").Value = 2026050.5
$ws.Range("I99").Value = 2633265.8
$ws.Range("J99").Value = 1999.6666
$ws.Range("K99").Value = 2633265.8
$ws.Range("L99").Value = 1999.6666
$ws.Range("M99").Value = -2631767.8
$ws.Range("N99").Value = -4995.6666
$ws.Range("H107").Value = 406.2
$ws.Range("I107").Value = 244.3
$ws.Range("J107").Value = 730
$ws.Range("K107").Value = 244.3
$ws.Range("L107").Value = 730
$ws.Range("M107").Value = 1675.7
$ws.Range("N107").Value = -4570
$ws.Range("H126").Value = 2026050.5
$ws.Range("I126").Value = 2633265.8
$ws.Range("J126").Value = 1999.6666
$ws.Range("K126").Value = 7899797.399999999
$ws.Range("L126").Value = 5998.9998
$ws.Range("M126").Value = -7897327.399999999
$ws.Range("N126").Value = -10938.9998

# --- Sheet 5: CUL ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H131").Value = 15387102
$ws.Range("I131").Value = 111111510
$ws.Range("J131").Value = 2821.6606
$ws.Range("K131").Value = 333334530
$ws.Range("L131").Value = 8464.981800000001
$ws.Range("M131").Value = -333329490
$ws.Range("N131").Value = -18544.9818

# --- Sheet 6: GSM ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H70").Value = 56252650
$ws.Range("I70").Value = 41669430
$ws.Range("J70").Value = 100002300
$ws.Range("K70").Value = 41669430
$ws.Range("L70").Value = 100002300
$ws.Range("M70").Value = -41669160
$ws.Range("N70").Value = -100002840
$ws.Range("H73").Value = 56252650
$ws.Range("I73").Value = 41669430
$ws.Range("J73").Value = 100002300
$ws.Range("K73").Value = 41669430
$ws.Range("L73").Value = 100002300
$ws.Range("M73").Value = -41668494
$ws.Range("N73").Value = -100004172
$ws.Range("H80").Value = 2926.9333
$ws.Range("J80").Value = 3614.1428
$ws.Range("L80").Value = 3614.1428
$ws.Range("N80").Value = -5610.1428
$ws.Range("H83").Value = 2926.9333
$ws.Range("J83").Value = 3614.1428
$ws.Range("L83").Value = 18070.714
$ws.Range("N83").Value = -28054.714
$ws.Range("H97").Value = 641.8
$ws.Range("I97").Value = 677.25
$ws.Range("J97").Value = 500
$ws.Range("K97").Value = 677.25
$ws.Range("L97").Value = 500
$ws.Range("M97").Value = -181.25
$ws.Range("N97").Value = -1492
$ws.Range("H122").Value = 1644.3334
$ws.Range("I122").Value = 1785.2858
$ws.Range("J122").Value = 1362.4286
$ws.Range("K122").Value = 5355.857400000001
$ws.Range("L122").Value = 4087.2858
$ws.Range("M122").Value = -2905.857400000001
$ws.Range("N122").Value = -8987.2858
$ws.Range("H132").Value = 1948.6888
$ws.Range("I132").Value = 1414.7307
$ws.Range("K132").Value = 4244.1921
$ws.Range("M132").Value = -1714.1921

# --- Sheet 7: LTW ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H7").Value = 1603.75
$ws.Range("I7").Value = 1410.9375
$ws.Range("K7").Value = 1410.9375
$ws.Range("M7").Value = -1298.9375
$ws.Range("H40").Value = 2523.3333
$ws.Range("I40").Value = 2135.9092
$ws.Range("K40").Value = 2135.9092
$ws.Range("M40").Value = -1999.9092
$ws.Range("H46").Value = 2349.8333
$ws.Range("I46").Value = 1033.3334
$ws.Range("K46").Value = 1033.3334
$ws.Range("M46").Value = -845.3334
$ws.Range("H55").Value = 218.32143
$ws.Range("I55").Value = 196.73334
$ws.Range("J55").Value = 243.23077
$ws.Range("K55").Value = 196.73334
$ws.Range("L55").Value = 243.23077
$ws.Range("M55").Value = -23.73334
$ws.Range("N55").Value = -589.23077
$ws.Range("H122").Value = 10496125
$ws.Range("I122").Value = 20240180
$ws.Range("J122").Value = 2526.3845
$ws.Range("K122").Value = 60720540
$ws.Range("L122").Value = 7579.1535
$ws.Range("M122").Value = -60718090
$ws.Range("N122").Value = -12479.1535
$ws.Range("H126").Value = 1603.75
$ws.Range("I126").Value = 1410.9375
$ws.Range("K126").Value = 4232.8125
$ws.Range("M126").Value = -1762.8125
$ws.Range("H136").Value = 1827.091
$ws.Range("I136").Value = 2057.6
$ws.Range("J136").Value = 1635
$ws.Range("K136").Value = 6172.799999999999
$ws.Range("L136").Value = 4905
$ws.Range("M136").Value = -3622.799999999999
$ws.Range("N136").Value = -10005

# --- Sheet 8: WVR ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H107").Value = 602.8570999999999
$ws.Range("I107").Value = 484
$ws.Range("K107").Value = 1452
$ws.Range("M107").Value = 468
$ws.Range("H122").Value = 63001070
$ws.Range("I122").Value = 78750940
$ws.Range("J122").Value = 1575
$ws.Range("K122").Value = 236252820
$ws.Range("L122").Value = 4725
$ws.Range("M122").Value = -236250370
$ws.Range("N122").Value = -9625
